$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Question 3 intro sentence: drop "at 20 seconds " before
#    "the balloon has a volume of 500 mL."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "at 20 seconds the balloon has a volume", $true, $false, $false, $false, $false,
    $true, 1, $false, "the balloon has a volume", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the blank paragraph that used to separate the intro
#    sentence from the "Under constant pressure..." paragraph.
#    (It is the paragraph right after the intro sentence, #11.)
# ------------------------------------------------------------------
$d.Paragraphs(11).Range.Delete()

# ------------------------------------------------------------------
# 3. Drop the leading tab before "Under constant pressure and
#    temperature: " so the paragraph starts directly with the text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "^tUnder constant pressure", $true, $false, $false, $false, $false,
    $true, 1, $false, "Under constant pressure", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Replace the two numbered-list questions (now paragraphs 12 & 13)
#    with two plain (non-numbered) paragraphs indented the same
#    amount the list body used to be.
# ------------------------------------------------------------------
$d.Paragraphs(13).Range.Delete()
$d.Paragraphs(12).Range.Delete()

$underPara = $d.Paragraphs(11)
$underPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs(12)
$newPara1.Format.LeftIndent = 36
$newPara1.Format.FirstLineIndent = 0
$newPara1.Range.Text = "What is the volume of the balloon after 30 seconds?"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(13)
$newPara2.Format.LeftIndent = 36
$newPara2.Format.FirstLineIndent = 0
